# Rows 2,4,5,6,7,8,9,10,11,12,13 are being permuted: each destination row
# receives the full original content of a different source row (row 3 and
# row 14 are untouched). Use far-away staging rows to snapshot every
# source row's full content (A:AY) before any destination is overwritten,
# since several rows are both a source and a destination (this is a
# permutation with cycles, not a simple 1:1 swap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (content that must end up in destination)
$srcForDest = @{
    2  = 10
    4  = 2
    5  = 7
    6  = 11
    7  = 13
    8  = 4
    9  = 12
    10 = 9
    11 = 6
    12 = 8
    13 = 5
}

$stagingOffset = 1000

# Step 1: snapshot every distinct source row into a staging row far below
# the real data, before any destination row gets overwritten.
$sources = @(10, 2, 7, 11, 13, 4, 12, 9, 6, 8, 5)
foreach ($r in $sources) {
    $srcRange = $ws.Range("A" + $r + ":AY" + $r)
    $stageRow = $stagingOffset + $r
    $dstRange = $ws.Range("A" + $stageRow + ":AY" + $stageRow)
    $srcRange.Copy($dstRange)
}

# Step 2: clear each destination row fully so no stray leftover cells
# remain from its previous contents (the source row may have fewer
# populated cells than the row being overwritten).
foreach ($d in $srcForDest.Keys) {
    $ws.Range("A" + $d + ":AY" + $d).ClearContents()
}

# Step 3: copy staged snapshots into their final destination rows.
foreach ($d in $srcForDest.Keys) {
    $s = $srcForDest[$d]
    $stageRow = $stagingOffset + $s
    $stageRange = $ws.Range("A" + $stageRow + ":AY" + $stageRow)
    $dstRange = $ws.Range("A" + $d + ":AY" + $d)
    $stageRange.Copy($dstRange)
}

# Step 4: clear out the staging rows so they don't linger in the sheet.
foreach ($r in $sources) {
    $stageRow = $stagingOffset + $r
    $ws.Range("A" + $stageRow + ":AY" + $stageRow).ClearContents()
}

Write-Host "row rotation complete"
